$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for the MSME Participation table (row 8)
$ws.Range("B8").Value = "Number of employees"
$ws.Range("C8").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D8").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B8:D8").Style = "title"

# Micro row (row 9)
$ws.Range("A9").Value = "Micro"

# Small row (row 10)
$ws.Range("A10").Value = "Small"
$ws.Range("C10").Value = "> USD 1,000"

# Medium row (row 11)
$ws.Range("A11").Value = "Medium"
$ws.Range("B11").Value = "<500"
$ws.Range("C11").Value = "< USD 5,000"

# Large row (row 12)
$ws.Range("A12").Value = "Large"
$ws.Range("B12").Value = ">500"
$ws.Range("C12").Value = "> USD 5,000"
